# JobPlanning.xlsx — AShot without inspirepak and jenkins propfile
#
# The sheet lists, per Component (col A), a series of planned Activities
# (col B) together with their ScheduledResource / MISWorkCenter /
# PlannedResource details (cols D, L, M) and PlannedQty (col G).
#
# This edit re-orders a handful of activity rows so the data is consistent:
#   * rows 4-5  ("- - 2) 591346  2p"): "-" (Press Approval) now comes
#     before "Imposition" (matching the order already used in rows 2-3).
#   * rows 6-8  ("1) 591345  2p"): the Sheet-fed Press / Die Cutting / Cut
#     activities are rotated up one slot.
#   * rows 16-17 ("2) 591346 2p Packed"): "Boxing" now comes before
#     "Pallet" (matching the order already used in rows 10-11).
#
# (Rows 12-14 keep the same activities; they only shift position inside
# the workbook's internal shared-string table, with no visible change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 4 & 5: swap "Imposition" <-> "-" / Press Approval details ---
$ws.Range("B4").Value = "-"
$ws.Range("D4").Value = "169-Press Approval Task "
$ws.Range("L4").Value = "Press Approval Task"
$ws.Range("M4").Value = "169-Press Approval Task "

$ws.Range("B5").Value = "Imposition"
$ws.Range("D5").Value = "134-Prepare files for CTP"
$ws.Range("L5").Value = "134-Prepare files for CTP"
$ws.Range("M5").Value = "134-Prepare files for CTP"

# --- Rows 6, 7 & 8: rotate Sheet-fed Press / Die Cutting / Cut up by one ---
# (PlannedQty in column G holds numbers-as-text like "10,050" in this sheet,
# so the column is kept on the Text number format to stop it turning into
# a real number when re-assigned.)
$ws.Range("B6").Value = "Die Cutting"
$ws.Range("D6").Value = "462-Bobst Letterpress"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "10,050"
$ws.Range("L6").Value = "462-Bobst Letterpress"
$ws.Range("M6").Value = "462-Bobst Letterpress`n463-Brausse Diecutter"

$ws.Range("B7").Value = "Cut"
$ws.Range("D7").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "10,000"
$ws.Range("L7").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("M7").Value = "406-45`" Polar 115ED Cutter`n404-45`" Polar 115EMC Cutter`n405-54`" Polar 137EMC Cutter`n402-45`" Polar 115EMC Cutter`n403-54`" Polar 137ED Cutter"

$ws.Range("B8").Value = "Sheet-fed Press F 4x0"
$ws.Range("D8").Value = "203-Heid 105 4-Color"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "10,196"
$ws.Range("L8").Value = "203-Heid 105 4-Color"
$ws.Range("M8").Value = "203-Heid 105 4-Color"

# --- Row 16 & 17: swap "Pallet" <-> "Boxing" details ---
$ws.Range("B16").Value = "Boxing"
$ws.Range("D16").Value = "481-PackSize Boxmaker"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "20"
$ws.Range("L16").Value = "481-PackSize Boxmaker"
$ws.Range("M16").Value = "481-PackSize Boxmaker"

$ws.Range("B17").Value = "Pallet"
$ws.Range("D17").Value = "602-Skid Wrap Operator"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "1"
$ws.Range("L17").Value = "602-Skid Wrap Operator"
$ws.Range("M17").Value = "602-Skid Wrap Operator"
